# Add a new "time_taken" metadata column (F) to the panel worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 — same look as the other header cells (copy format from E1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$timestamps = @(
  "2021-10-05 10:50:53.612013",
  "2021-10-05 10:50:53.612024",
  "2021-10-05 10:50:53.612027",
  "2021-10-05 10:50:53.612030",
  "2021-10-05 10:50:53.612033",
  "2021-10-05 10:50:53.612036",
  "2021-10-05 10:50:53.612038",
  "2021-10-05 10:50:53.612041",
  "2021-10-05 10:50:53.612043",
  "2021-10-05 10:50:53.612046",
  "2021-10-05 10:50:53.612048",
  "2021-10-05 10:50:53.612051",
  "2021-10-05 10:50:53.612053",
  "2021-10-05 10:50:53.612056",
  "2021-10-05 10:50:53.612058",
  "2021-10-05 10:50:53.612061",
  "2021-10-05 10:50:53.612064",
  "2021-10-05 10:50:53.612066",
  "2021-10-05 10:50:53.612069",
  "2021-10-05 10:50:53.612071",
  "2021-10-05 10:50:53.612074",
  "2021-10-05 10:50:53.612076",
  "2021-10-05 10:50:53.612079",
  "2021-10-05 10:50:53.612081",
  "2021-10-05 10:50:53.612084",
  "2021-10-05 10:50:53.612087",
  "2021-10-05 10:50:53.612089",
  "2021-10-05 10:50:53.612092",
  "2021-10-05 10:50:53.612094",
  "2021-10-05 10:50:53.612096",
  "2021-10-05 10:50:53.612099",
  "2021-10-05 10:50:53.612101"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

$ws.Range("A1").Select() | Out-Null
